$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from serial date 45174 to 45175 for data rows 2-12
$ws.Range("C2:C12").Value = 45175
